# testActors.xlsx: rename the "movie ID" header (A1) to "movie", and move
# the active selection from E11 to E9 (matches the authored diff; the
# sharedStrings reshuffle in the diff is just a side-effect of Excel
# rewriting the string table once "movie ID" stops being referenced).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "movie"

$null = $ws.Range("E9").Select()
